$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.309671333333333
$ws.Range("N2").Value = 3.929014
$ws.Range("O2").Value = 0.05806924226264097
$ws.Range("P2").Value = 0.05806924226264098
$ws.Range("Q2").Value = 0.195004822848
$ws.Range("R2").Value = 1.755043405632
$ws.Range("S2").Value = 0.05806924226264097
$ws.Range("T2").Value = 0.05806924226264098

# Row 3
$ws.Range("O3").Value = 0.3245116581089107
$ws.Range("P3").Value = 0.3245116581089107
$ws.Range("R3").Value = 9.807809150303999
$ws.Range("S3").Value = 0.3245116581089107
$ws.Range("T3").Value = 0.3245116581089107

# Row 4
$ws.Range("M4").Value = 4.657910333333334
$ws.Range("N4").Value = 13.973731
$ws.Range("O4").Value = 0.2065261082683789
$ws.Range("P4").Value = 0.2065261082683789
$ws.Range("Q4").Value = 0.6935442169920001
$ws.Range("R4").Value = 6.241897952928
$ws.Range("S4").Value = 0.2065261082683789
$ws.Range("T4").Value = 0.2065261082683789

# Row 5
$ws.Range("M5").Value = 9.267122333333333
$ws.Range("N5").Value = 27.801367
$ws.Range("O5").Value = 0.4108929913600695
$ws.Range("P5").Value = 0.4108929913600695
$ws.Range("Q5").Value = 1.379837446944
$ws.Range("R5").Value = 12.418537022496
$ws.Range("S5").Value = 0.4108929913600695
$ws.Range("T5").Value = 0.4108929913600695
